$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 101, pushing existing rows 101-150 down to 102-151
$ws.Rows("101:101").Insert()

# Populate the newly inserted row 101 with the new report entry
$ws.Range("A101").Value = 10
$ws.Range("B101").Value = "Vega Modelo de Temuco"
$ws.Range("C101").Value = "La Araucanía"
$ws.Range("D101").Value = 44523
$ws.Range("E101").Value = 9
$ws.Range("F101").Value = 100112052
$ws.Range("G101").Value = "Albahaca"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 30
$ws.Range("K101").Value = 3500
$ws.Range("L101").Value = 3500
$ws.Range("M101").Value = 3500
$ws.Range("N101").Value = "$/paquete"
$ws.Range("O101").Value = "Región del Maule"
$ws.Range("P101").Value = 3500
$ws.Range("Q101").Value = 1
$ws.Range("R101").Value = "Hortaliza"
